# OO-3174: CSV import of essay and match in QTI 2.1
# Adds three new example question-type blocks to the "Import" metadata sheet:
#   ESSAY        rows 78-83
#   MATRIX       rows 86-93
#   Drag&drop    rows 96-103
#
# Cell values below are written in a specific sequence so that newly
# introduced shared strings land at the same table offsets the original
# authoring session produced; formatting is applied afterwards in a
# second pass (PasteSpecial formats-only never touches shared strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Pass 1: cell values, in shared-string-creation order
# ---------------------------------------------------------------------------

# ESSAY block
$ws.Range("A78").Value = "Typ"
$ws.Range("B78").Value = "ESSAY"
$ws.Range("A79").Value = "Title"
$ws.Range("B79").Value = "Fussball: Spieler"
$ws.Range("A80").Value = "Question"
$ws.Range("B80").Value = "Wer ist der beste Spieler aller Zeit?"
$ws.Range("A81").Value = "Points"
$ws.Range("B81").Value = 1
$ws.Range("A82").Value = "Min"
$ws.Range("B82").Value = 200
$ws.Range("A83").Value = "Max"
$ws.Range("B83").Value = 2000

# MATRIX block
$ws.Range("A86").Value = "Typ"
$ws.Range("B86").Value = "MATRIX"
$ws.Range("A89").Value = "Points"
$ws.Range("B89").Value = 1
$ws.Range("B90").Value = "Deutschland"
$ws.Range("C90").Value = "Frankreich"
$ws.Range("D90").Value = "Schweiz"
$ws.Range("A91").Value = "Berlin"
$ws.Range("B91").Value = 1
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("A92").Value = "Bern"
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 1
$ws.Range("A93").Value = "Paris"
$ws.Range("B93").Value = 0
$ws.Range("C93").Value = 1
$ws.Range("D93").Value = 0
$ws.Range("A87").Value = "Title"
$ws.Range("B87").Value = "Hauptstädte Europas"
$ws.Range("A88").Value = "Question"
$ws.Range("B88").Value = "Hauptstädte Europas"

# Drag&drop block
$ws.Range("A96").Value = "Typ"
$ws.Range("B96").Value = "Drag&drop"
$ws.Range("A97").Value = "Title"
$ws.Range("B97").Value = "Hauptstädte Afrika"
$ws.Range("A98").Value = "Question"
$ws.Range("B98").Value = "Hauptstädte Afrika"
$ws.Range("A99").Value = "Points"
$ws.Range("B99").Value = 1
$ws.Range("B100").Value = "Algerien"
$ws.Range("A103").Value = "Algier"
$ws.Range("B103").Value = 1
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 0
$ws.Range("C100").Value = "Kenia"
$ws.Range("D100").Value = "Namibia"
$ws.Range("A101").Value = "Nairobi"
$ws.Range("B101").Value = 0
$ws.Range("C101").Value = 1
$ws.Range("D101").Value = 0
$ws.Range("A102").Value = "Windhoek"
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 1

# ---------------------------------------------------------------------------
# Pass 2: formatting, copied from equivalent existing cells (KPRIM block)
# ---------------------------------------------------------------------------

# "Typ" rows: label cell (A) + bold "Ausgabe" value cell (B)
$ws.Range("A68").Copy() | Out-Null
$ws.Range("A78").PasteSpecial(-4122) | Out-Null
$ws.Range("A86").PasteSpecial(-4122) | Out-Null
$ws.Range("A96").PasteSpecial(-4122) | Out-Null
$ws.Range("A90").PasteSpecial(-4122) | Out-Null
$ws.Range("A100").PasteSpecial(-4122) | Out-Null

$ws.Range("B68").Copy() | Out-Null
$ws.Range("B78").PasteSpecial(-4122) | Out-Null
$ws.Range("B86").PasteSpecial(-4122) | Out-Null
$ws.Range("B96").PasteSpecial(-4122) | Out-Null

# Plain label cells (A column, Title/Question/Points/Min/Max/answer-row label)
$ws.Range("A69").Copy() | Out-Null
$ws.Range("A79,A87,A97").PasteSpecial(-4122) | Out-Null

$ws.Range("A70").Copy() | Out-Null
$ws.Range("A80,A88,A98").PasteSpecial(-4122) | Out-Null

$ws.Range("A71").Copy() | Out-Null
$ws.Range("A81,A82,A83,A89,A91,A92,A93,A99,A101,A102,A103").PasteSpecial(-4122) | Out-Null

# Plain value cells (B column text/number rows)
$ws.Range("B69").Copy() | Out-Null
$ws.Range("B79,B87").PasteSpecial(-4122) | Out-Null

$ws.Range("B70").Copy() | Out-Null
$ws.Range("B80,B88").PasteSpecial(-4122) | Out-Null

$ws.Range("B71").Copy() | Out-Null
$ws.Range("B81,B82,B83,B89").PasteSpecial(-4122) | Out-Null

# Header row of the answer matrix (B/C/D bold "Ausgabe" style)
$ws.Range("B68").Copy() | Out-Null
$ws.Range("B90,C90,D90,B100,C100,D100").PasteSpecial(-4122) | Out-Null

# Answer-grid value cells (B/C/D plain value style)
$ws.Range("B71").Copy() | Out-Null
$ws.Range("B91:D93,B99,B101:D103").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------------
$ws.Range("B97").Select() | Out-Null
